# Update (Analyze PO & Forecast)
#
# Rewrites the 16-week forecast table on "Forecast Comparison" with a
# refreshed pull (week-start dates shifted back 4 weeks, new MyForecast
# values, and is_holiday_week cleared from boolean FALSE to a blank/zero
# numeric), then refreshes the dependent roll-up figures on "Summary".

$wb = $excel.ActiveWorkbook
$wsFC = $wb.Worksheets("Forecast Comparison")
$wsSum = $wb.Worksheets("Summary")

# Helper: write $text into $cellRef as literal TEXT (never let Excel's
# autocomplete reinterpret an ISO date / numeric-looking string as a real
# date serial or number). We stage the value in a scratch cell formatted
# as Text, copy it, and paste-special *values only* so no number format
# gets carried onto the destination cell (keeps styles untouched, matching
# how the sheet originally had no explicit per-cell formatting).
function Set-TextValue {
    param($ws, $cellRef, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# New Week_Start_Date (col B) and MyForecast (col D) values for rows 2-17.
$weekStart = @(
    "2024-12-29", "2025-01-05", "2025-01-12", "2025-01-19",
    "2025-01-26", "2025-02-02", "2025-02-09", "2025-02-16",
    "2025-02-23", "2025-03-02", "2025-03-09", "2025-03-16",
    "2025-03-23", "2025-03-30", "2025-04-06", "2025-04-13"
)
$myForecast = @(63, 61, 50, 54, 62, 55, 57, 60, 60, 58, 59, 61, 61, 58, 58, 58)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $bRef = "B" + $row

    Set-TextValue $wsFC $bRef $weekStart[$i]
    $wsFC.Range("D" + $row).Value = $myForecast[$i]
    # is_holiday_week: was boolean FALSE, now comes through as a blank
    # numeric (0) from the refreshed pull.
    $wsFC.Range("J" + $row).Value = 0
}

# Dependent roll-up figures on the Summary sheet.
Set-TextValue $wsSum "B9"  "935"
Set-TextValue $wsSum "B10" "462"
Set-TextValue $wsSum "B11" "228"
Set-TextValue $wsSum "B12" "63"
Set-TextValue $wsSum "B13" "2024-12-29"
Set-TextValue $wsSum "B14" "50"
Set-TextValue $wsSum "B15" "2025-01-12"
